$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5187906666666667
$ws.Range("H2").Value = 1.556372
$ws.Range("I2").Value = 0.5259328345599914
$ws.Range("J2").Value = 0.5259328345599914
$ws.Range("M2").Value = 32.21373866666666
$ws.Range("N2").Value = 96.641216
$ws.Range("O2").Value = 0.6812298485843321
$ws.Range("P2").Value = 0.7117693664123
$ws.Range("Q2").Value = 16.71218695870578
$ws.Range("R2").Value = 150.409682628352
$ws.Range("S2").Value = 0.3582811452528316
$ws.Range("T2").Value = 0.3743428804301901

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5187906666666667
$ws.Range("H3").Value = 1.556372
$ws.Range("I3").Value = 0.5259328345599914
$ws.Range("J3").Value = 0.5259328345599914
$ws.Range("O3").Value = 0.1770389772624213
$ws.Range("P3").Value = 0.184975630381169
$ws.Range("Q3").Value = 4.343186801248001
$ws.Range("R3").Value = 39.08868121123201
$ws.Range("S3").Value = 0.09311061113922708
$ws.Range("T3").Value = 0.09728475761088949

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5187906666666667
$ws.Range("H4").Value = 1.556372
$ws.Range("I4").Value = 0.5259328345599914
$ws.Range("J4").Value = 0.5259328345599914
$ws.Range("M4").Value = 0.5484013333333334
$ws.Range("N4").Value = 1.645204
$ws.Range("O4").Value = 0.01159714372603029
$ws.Range("P4").Value = 0.01211704340205096
$ws.Range("Q4").Value = 0.284505493320889
$ws.Range("R4").Value = 2.560549439888
$ws.Range("S4").Value = 0.006099318672630729
$ws.Range("T4").Value = 0.006372750982927103

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5187906666666667
$ws.Range("H5").Value = 1.556372
$ws.Range("I5").Value = 0.5259328345599914
$ws.Range("J5").Value = 0.5259328345599914
$ws.Range("M5").Value = 6.086836
$ws.Range("N5").Value = 12.173672
$ws.Range("O5").Value = 0.1287194389184112
$ws.Range("P5").Value = 0.08965995219214913
$ws.Range("Q5").Value = 3.157793706330667
$ws.Range("R5").Value = 18.946762237984
$ws.Range("S5").Value = 0.06769777937333168
$ws.Range("T5").Value = 0.0471551128029303

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5187906666666667
$ws.Range("H6").Value = 1.556372
$ws.Range("I6").Value = 0.5259328345599914
$ws.Range("J6").Value = 0.5259328345599914
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06689266666666667
$ws.Range("N6").Value = 0.200678
$ws.Range("O6").Value = 0.001414591508805173
$ws.Range("P6").Value = 0.001478007612330618
$ws.Range("Q6").Value = 0.03470329113511111
$ws.Range("R6").Value = 0.312329620216
$ws.Range("S6").Value = 0.0007439801219703996
$ws.Range("T6").Value = 0.0007773327330542869

# Row 7
$ws.Range("G7").Value = 0.4676293333333333
$ws.Range("H7").Value = 1.402888
$ws.Range("I7").Value = 0.4740671654400085
$ws.Range("J7").Value = 0.4740671654400086
$ws.Range("M7").Value = 32.21373866666666
$ws.Range("N7").Value = 96.641216
$ws.Range("O7").Value = 0.6812298485843321
$ws.Range("P7").Value = 0.7117693664123
$ws.Range("Q7").Value = 15.06408913686755
$ws.Range("R7").Value = 135.576802231808
$ws.Range("S7").Value = 0.3229487033315006
$ws.Range("T7").Value = 0.3374264859821099

# Row 8
$ws.Range("G8").Value = 0.4676293333333333
$ws.Range("H8").Value = 1.402888
$ws.Range("I8").Value = 0.4740671654400085
$ws.Range("J8").Value = 0.4740671654400086
$ws.Range("O8").Value = 0.1770389772624213
$ws.Range("P8").Value = 0.184975630381169
$ws.Range("Q8").Value = 3.914876806592
$ws.Range("R8").Value = 35.233891259328
$ws.Range("S8").Value = 0.08392836612319417
$ws.Range("T8").Value = 0.08769087277027955

# Row 9
$ws.Range("G9").Value = 0.4676293333333333
$ws.Range("H9").Value = 1.402888
$ws.Range("I9").Value = 0.4740671654400085
$ws.Range("J9").Value = 0.4740671654400086
$ws.Range("M9").Value = 0.5484013333333334
$ws.Range("N9").Value = 1.645204
$ws.Range("O9").Value = 0.01159714372603029
$ws.Range("P9").Value = 0.01211704340205096
$ws.Range("Q9").Value = 0.2564485499057778
$ws.Range("R9").Value = 2.308036949152
$ws.Range("S9").Value = 0.005497825053399557
$ws.Range("T9").Value = 0.005744292419123858

# Row 10
$ws.Range("G10").Value = 0.4676293333333333
$ws.Range("H10").Value = 1.402888
$ws.Range("I10").Value = 0.4740671654400085
$ws.Range("J10").Value = 0.4740671654400086
$ws.Range("M10").Value = 6.086836
$ws.Range("N10").Value = 12.173672
$ws.Range("O10").Value = 0.1287194389184112
$ws.Range("P10").Value = 0.08965995219214913
$ws.Range("Q10").Value = 2.846383060789333
$ws.Range("R10").Value = 17.078298364736
$ws.Range("S10").Value = 0.06102165954507953
$ws.Range("T10").Value = 0.04250483938921883

# Row 11
$ws.Range("G11").Value = 0.4676293333333333
$ws.Range("H11").Value = 1.402888
$ws.Range("I11").Value = 0.4740671654400085
$ws.Range("J11").Value = 0.4740671654400086
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.06689266666666667
$ws.Range("N11").Value = 0.200678
$ws.Range("O11").Value = 0.001414591508805173
$ws.Range("P11").Value = 0.001478007612330618
$ws.Range("Q11").Value = 0.03128097311822222
$ws.Range("R11").Value = 0.281528758064
$ws.Range("S11").Value = 0.0007006748792763313
$ws.Range("T11").Value = 0.0007006748792763313
